$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'76.110.25"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = "'2.917.43"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('D5').Value = "'202.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.35%  '
$ws.Range('D6').Value = "'598.31"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('D10').Value = "'2.914.95"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.47%  '
$ws.Range('E11').Value = '  +17.91%  '
$ws.Range('D12').Value = "'0.162"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = "'3.452.60"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').Value = "'75.956.88"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'27.95"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.34%  '
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = "'2.906.18"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').Value = "'12.91"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.17%  '
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').Value = "'373.68"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.74%  '
$ws.Range('D22').Value = "'2.32"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.38%  '
$ws.Range('E23').Value = '  +6.22%  '
$ws.Range('D24').Value = "'71.39"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = "'3.053.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('D28').Value = "'9.70"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').Value = "'504.79"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = "'EthereumClassic"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'20.24"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('B37').Value = "'Monero"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'163.82"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('E38').Value = '  +1.26%  '
$ws.Range('D39').Value = "'0.106"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +24.21%  '
$ws.Range('E40').Value = '  -4.56%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'181.33"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.36%  '
$ws.Range('D43').Value = "'0.359"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.08%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = "'1.66"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = "'39.99"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = "'1.19"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').Value = "'0.572"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('B51').Value = "'Mantle"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'0.654"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.60%  '
